# Update column AA ("average_county_temperature") with refreshed NOAA
# temperature readings for the affected facility rows.
#
# The underlying source data was re-pulled from NOAA; the rows below map to
# the facility groups whose average_county_temperature changed as a result
# (facility 1000606, 1000615, 1001985, 1002283, 1002285, 1004369, 1005361).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    @{ Rows = 2..5;     Value = 1.925925925925943 },
    @{ Rows = 50..89;   Value = 1.925925925925943 },
    @{ Rows = 90..106;  Value = -1.226851851851833 },
    @{ Rows = 119..130; Value = 1.925925925925943 },
    @{ Rows = 175..206; Value = 13.17361111111111 }
)

foreach ($group in $ranges) {
    foreach ($row in $group.Rows) {
        $ws.Range("AA$row").Value = $group.Value
    }
}
